$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Raiden / 雷专) : fix a stray extra zero in fixed-HP value ---
$ws.Range("Q2").Value = 4780

# --- Row 3 (Bennett / 班尼特) : switch artifact set to Noblesse Oblige (宗室),
#     update skill levels + stats. Weapon text (天空剑) is set last, below,
#     so new shared-string slots line up with how the workbook was authored. ---
$ws.Range("E3").Value = "9,9,9"
$ws.Range("H3").Value = "宗室"
$ws.Range("I3").Value = "宗室"
$ws.Range("K3").Value = 311
$ws.Range("Q3").Value = 4780
$ws.Range("R3").Value = 93.2

# --- Row 4 (new): Kamisato Ayaka (神里绫华) built with Mistsplitter Reforged
#     (雾切) and Blizzard Strayer (冰套) ---
$ws.Range("B4").Value = "神里绫华"
$ws.Range("A4").Value = "绫华"
$ws.Range("C4").Value = "90+"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "9,9,10"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "冰套"
$ws.Range("I4").Value = "冰套"
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 311
$ws.Range("L4").Value = 46.6
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 4780
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 40
$ws.Range("U4").Value = 140
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0
$ws.Range("AO4").Value = 0
$ws.Range("AP4").Value = 0

# Weapon text for row 4 (Mistsplitter, 雾切) then row 3 (Skyward Blade, 天空剑)
# last, matching the order new strings were appended in the authored workbook.
$ws.Range("F4").Value = "雾切"
$ws.Range("F3").Value = "天空剑"

# --- View state: move the active selection to F4 and set normal-view zoom ---
[void]$ws.Select()
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("F4").Select()
